$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '59.375.23'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  +0.27%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '2.637.32'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  +1.32%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  +0.01%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '536.61'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  -0.82%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '145.28'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  +2.70%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  -0.04%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  +1.13%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  +9.47%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'" + '0.101'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  -1.20%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'" + '0.338'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  +0.83%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  +0.26%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'" + '3.103.53'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  +1.42%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'" + '59.302.75'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  +0.25%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'" + '21.27'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  +3.39%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'" + '2.640.22'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  +1.58%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  +3.31%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'" + '338.32'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  -0.81%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'" + '10.29'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  -2.37%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  +0.10%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'" + '66.21'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  -2.12%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'" + '0.416'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  +1.87%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  +0.01%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'" + '0.989'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  -1.00%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  +1.47%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'" + '0.0₃0752'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  -0.04%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'" + '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  -0.05%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'" + '1.65'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Value = "'" + '5.90'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  +1.40%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'" + '18.83'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  +0.55%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'" + '151.05'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  +0.89%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'" + '4.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  +0.39%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  +2.13%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'" + '0.843'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  +2.32%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'" + '0.839'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  +0.60%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  -1.07%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  +1.24%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'" + '285.17'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  +3.80%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  -0.02%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'" + '0.600'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  +0.50%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  +0.07%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  +2.95%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'" + '19.13'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  +2.59%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'" + '0.0943'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  -1.45%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  +1.63%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'" + '1.960.82'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  +0.32%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = "'" + 'RenderToken'
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = "'" + 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = "'" + '4.55'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  +0.93%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = "'" + 'InjectiveProtocol'
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = "'" + 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = "'" + '18.39'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  -0.87%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'" + '111.44'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  +0.15%  '
$ws.Range('E51').Style = 'Normal'
